# ProductsList.xlsx edit:
#  - populate newly-added product_sub_category ("None") cells C25:C33
#  - move the sheet's viewport/selection to rows 15+, selecting C25:C33
#  - refresh the "Y" conditional-formatting rule on F2:F33 (dxf cleanup)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fill in column C ("product_sub_category") for rows 25-33 with "None",
#    matching the same value already used for all the other rows above.
$ws.Range("C25:C33").Value = "None"

# 2. Select C25:C33 and scroll the window so row 15 is at the top
#    (mirrors <sheetView topLeftCell="A15"> / <selection activeCell="C25" sqref="C25:C33"/>)
$ws.Activate() | Out-Null
$ws.Range("C25:C33").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 1

# 3. Re-apply the "equals Y" conditional format rule on F2:F33 (bold red
#    font on a light accent fill) - this is what the author touched in the
#    Conditional Formatting manager, which made Excel drop the now-unused
#    dxf records left behind from earlier edits of this rule.
$fConds = $ws.Range("F2:F33").FormatConditions
$rule = $fConds.Item(1)
$rule.Modify(1, 3, '"Y"')
$rule.Font.Bold = $true
$rule.Font.Italic = $false
$rule.Font.Color = 255
$rule.Interior.Color = 13431551

Write-Host "ProductsList.xlsx updated"
